$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 61 currently hold the value 7310 and need to
# become 7293 (values for the remaining rows are already 7293 and stay
# untouched).
$ws.Range("C2:C61").Value = 7293
